$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'281.56"
$ws.Range("G2").Value = "'6"
$ws.Range("D3").Value = "'20.64"
$ws.Range("G3").Value = "'6"
$ws.Range("D4").Value = "'6.212"
$ws.Range("G4").Value = "'6"
$ws.Range("D5").Value = "'0.06158"
$ws.Range("G5").Value = "'6"
$ws.Range("D6").Value = "'3.581"
$ws.Range("G6").Value = "'6"
$ws.Range("D7").Value = "'6.575"
$ws.Range("G7").Value = "'6"
$ws.Range("D8").Value = "'1.501"
$ws.Range("G8").Value = "'6"
$ws.Range("D9").Value = "'0.8193"
$ws.Range("G9").Value = "'6"
$ws.Range("D10").Value = "'0.01387"
$ws.Range("G10").Value = "'6"
$ws.Range("D11").Value = "'0.1645"
$ws.Range("G11").Value = "'6"
$ws.Range("D12").Value = "'0.08422"
$ws.Range("G12").Value = "'6"
$ws.Range("D13").Value = "'0.03518"
$ws.Range("G13").Value = "'6"
$ws.Range("D14").Value = "'0.03207"
$ws.Range("G14").Value = "'6"
$ws.Range("D15").Value = "'0.09134"
$ws.Range("G15").Value = "'6"
$ws.Range("D16").Value = "'3.704"
$ws.Range("G16").Value = "'6"
$ws.Range("D17").Value = "'0.001640"
$ws.Range("G17").Value = "'6"
$ws.Range("D18").Value = "'0.04728"
$ws.Range("G18").Value = "'6"
$ws.Range("D19").Value = "'0.006537"
$ws.Range("G19").Value = "'6"
$ws.Range("D20").Value = "'0.006163"
$ws.Range("G20").Value = "'6"
$ws.Range("D21").Value = "'0.001070"
$ws.Range("G21").Value = "'6"
$ws.Range("G22").Value = "'6"
$ws.Range("D23").Value = "'3.781"
$ws.Range("G23").Value = "'6"
$ws.Range("D24").Value = "'2.323"
$ws.Range("G24").Value = "'6"
$ws.Range("G25").Value = "'6"
$ws.Range("G26").Value = "'6"
$ws.Range("G27").Value = "'6"
$ws.Range("B28").Value = "Spectre.aiUtilityToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("D28").Value = "--"
$ws.Range("E28").Value = "27Spectre.aiUtilityTokenSXUT"
$ws.Range("G28").Value = "'6"
$ws.Range("B29").Value = "LegolasExchange"
$ws.Range("C29").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("E29").Value = "28LegolasExchangeLGO"
$ws.Range("G29").Value = "'6"
$ws.Range("B30").Value = "BitZToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("E30").Value = "29BitZTokenBZ"
$ws.Range("G30").Value = "'6"
$ws.Range("B31").Value = "Birake"
$ws.Range("C31").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("E31").Value = "30BirakeBIR"
$ws.Range("G31").Value = "'6"
$ws.Range("B32").Value = "ZBToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("E32").Value = "31ZBTokenZB"
$ws.Range("G32").Value = "'6"
$ws.Range("B33").Value = "NashExchange"
$ws.Range("C33").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("E33").Value = "32NashExchangeNEX"
$ws.Range("G33").Value = "'6"
$ws.Range("B34").Value = "CenX"
$ws.Range("C34").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("E34").Value = "33CenXCENX"
$ws.Range("G34").Value = "'6"
$ws.Range("B35").Value = "BNIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("E35").Value = "34BNIXTokenBNIX"
$ws.Range("G35").Value = "'6"
$ws.Range("B36").Value = "UpBots"
$ws.Range("C36").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("E36").Value = "35UpBotsUBXT"
$ws.Range("G36").Value = "'6"
$ws.Range("G37").Value = "'6"
$ws.Range("G38").Value = "'6"
$ws.Range("G39").Value = "'6"
$ws.Range("D40").Value = "'0.04697"
$ws.Range("G40").Value = "'6"
$ws.Range("D41").Value = "'0.007189"
$ws.Range("G41").Value = "'6"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1099"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'6"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003303"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'6"
$ws.Range("D44").Value = "'0.01107"
$ws.Range("G44").Value = "'6"
$ws.Range("D45").Value = "'0.00006536"
$ws.Range("G45").Value = "'6"
$ws.Range("G46").Value = "'6"
$ws.Range("D47").Value = "'1.001"
$ws.Range("G47").Value = "'6"
$ws.Range("D48").Value = "'0.002849"
$ws.Range("G48").Value = "'6"
$ws.Range("D49").Value = "'0.00001902"
$ws.Range("G49").Value = "'6"
$ws.Range("G50").Value = "'6"
$ws.Range("G51").Value = "'6"
